# Slide 8 ("JSX vs HTML"), Content Placeholder 2: append a new bullet
# paragraph made of two runs - "In JSX, we have to pass single React
# element to " followed by "return statement" - after the existing
# "JSX can reference JS variables" bullet.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(8)
$sh = $s.Shapes.Item("Content Placeholder 2")
$tr = $sh.TextFrame.TextRange

# Start a new paragraph at the end of the existing text, then type the
# first run of the new bullet.
$run1 = $tr.InsertAfter("`rIn JSX, we have to pass single React element to ")

# Continue typing immediately after it with the second run of the bullet.
$run2 = $run1.InsertAfter("return statement")
